$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.083.06"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "2.551.74"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'581.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.45%  "
$ws.Range("D6").Value = "'147.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.11%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +0.53%  "
$ws.Range("E10").Value = "  -1.80%  "
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("D13").Value = "'27.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.29%  "
$ws.Range("D14").Value = "3.007.08"
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("D15").Value = "62.971.34"
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("E16").Value = "  +1.63%  "
$ws.Range("D17").Value = "2.547.13"
$ws.Range("E17").Value = "  -2.62%  "
$ws.Range("E18").Value = "  -2.09%  "
$ws.Range("D19").Value = "'339.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("E21").Value = "  -0.41%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").Value = "'65.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.55%  "
$ws.Range("D24").Value = "2.675.76"
$ws.Range("E24").Value = "  +0.42%  "
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("E26").Value = "  +1.90%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").Value = "'1.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.10%  "
$ws.Range("D29").Value = "'8.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.84%  "
$ws.Range("D30").Value = "'7.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.25%  "
$ws.Range("D31").Value = "'1.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.13%  "
$ws.Range("D32").Value = "0.0₃0819"
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("D33").Value = "'178.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("D34").Value = "'423.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.28%  "
$ws.Range("E35").Value = "  -1.19%  "
$ws.Range("E36").Value = "  -1.16%  "
$ws.Range("E37").Value = "  +0.42%  "
$ws.Range("D39").Value = "'4.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.74%  "
$ws.Range("E40").Value = "  -0.95%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  +0.68%  "
$ws.Range("D43").Value = "'150.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.69%  "
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").Value = "'20.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.59%  "
$ws.Range("E46").Value = "  +3.39%  "
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("D49").Value = "'0.0240"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.05%  "
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("E51").Value = "  -4.71%  "
